# Append 12 more observations (months 204-215) to the normalized CPI series,
# continuing the existing pattern in columns A (index) and B (normalized value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Index (column A) / normalized value (column B) pairs to append.
$newData = @(
    @(204, 0.3941088739746457),
    @(205, 0.3536031455494542),
    @(206, 0.3444337914136572),
    @(207, 0.3999857959589503),
    @(208, 0.5186428038777031),
    @(209, 0.6141320062700695),
    @(210, 0.3778097368701395),
    @(211, 0.4836316181953765),
    @(212, 0.3944817300521998),
    @(213, 0.2456588899541919),
    @(214, 0.3447853414296367),
    @(215, 0.444018323212954)
)

# Use the last existing data row (205) as the style template for column A,
# so the new index cells keep the same formatting as the rest of the series.
$styleTemplate = $ws.Range("A205")
$styleTemplate.Copy()

$lastRow = 205
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $lastRow + 1 + $i
    $indexValue = $newData[$i][0]
    $normValue = $newData[$i][1]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.PasteSpecial(-4122)  # xlPasteFormats - replicate column A's style
    $aCell.Value = $indexValue

    $ws.Cells.Item($row, 2).Value = $normValue
}

$excel.CutCopyMode = 0
